# Commit: "new sim results and new calculation"
#
# 1) annualised_return / mean_period_return / sharpe_annualized get refreshed
#    simulation numbers.
# 2) A new "sharpe_period" sheet is inserted right before "VaR", and "VaR"
#    is recalculated with a brand new matrix of values. (Implementation
#    note: this is done by renaming the existing "VaR" sheet to
#    "sharpe_period" in place - picking up its header/label layout for free
#    - and adding a fresh sheet named "VaR" right after it.)

$wb = $excel.ActiveWorkbook
$cols = @("B","C","D","E","F","G","H","I","J","K")

function Set-MatrixValues {
    param($ws, $rowValues)
    foreach ($r in $rowValues.Keys) {
        $vals = $rowValues[$r]
        for ($i = 0; $i -lt $cols.Count; $i++) {
            $ws.Range("$($cols[$i])$r").Value = $vals[$i]
        }
    }
}

# ---------------------------------------------------------------------------
# 1) annualised_return - updated values
# ---------------------------------------------------------------------------
$wsAnn = $wb.Worksheets.Item("annualised_return")
Set-MatrixValues $wsAnn @{
    2  = @(1, 0, 0, 0, 0.0002, 0, 0, 0, 0, 0)
    3  = @(0, 1, 0, 0, 0, 0, 0, 0, 0, 0.5263)
    4  = @(0, 0, 1, 1, 0, 0, 0, 0, 0, 0)
    5  = @(0, 0, 1, 1, 0, 0, 0, 0, 0, 0)
    6  = @(0.0002, 0, 0, 0, 1, 0, 0, 0, 0, 0)
    7  = @(0, 0, 0, 0, 0, 1, 0, 0, 0, 0.4747)
    8  = @(0, 0, 0, 0, 0, 0, 1, 1, 1, 0)
    9  = @(0, 0, 0, 0, 0, 0, 1, 1, 1, 0)
    10 = @(0, 0, 0, 0, 0, 0, 1, 1, 1, 0)
    11 = @(0, 0.5263, 0, 0, 0, 0.4747, 0, 0, 0, 1)
}

# ---------------------------------------------------------------------------
# 2) mean_period_return - updated values
# ---------------------------------------------------------------------------
$wsMean = $wb.Worksheets.Item("mean_period_return")
Set-MatrixValues $wsMean @{
    2  = @(1, 0, 0, 0, 0.0009, 0, 0, 0, 0, 0)
    3  = @(0, 1, 0, 0, 0, 0, 0, 0, 0, 0.0239)
    4  = @(0, 0, 1, 1, 0, 0, 0, 0, 0, 0)
    5  = @(0, 0, 1, 1, 0, 0, 0, 0, 0, 0)
    6  = @(0.0009, 0, 0, 0, 1, 0, 0, 0, 0, 0)
    7  = @(0, 0, 0, 0, 0, 1, 0, 0, 0, 1)
    8  = @(0, 0, 0, 0, 0, 0, 1, 1, 1, 0)
    9  = @(0, 0, 0, 0, 0, 0, 1, 1, 1, 0)
    10 = @(0, 0, 0, 0, 0, 0, 1, 1, 1, 0)
    11 = @(0, 0.0239, 0, 0, 0, 1, 0, 0, 0, 1)
}

# ---------------------------------------------------------------------------
# 3) sharpe_annualized - updated values
# ---------------------------------------------------------------------------
$wsSharpeAnn = $wb.Worksheets.Item("sharpe_annualized")
Set-MatrixValues $wsSharpeAnn @{
    2  = @(1, 1, 1, 1, 0, 0, 0, 0, 0, 0)
    3  = @(1, 1, 0.4748, 0.4741, 0, 0, 0, 0, 0, 0)
    4  = @(1, 0.4748, 1, 1, 0, 0, 0, 0, 0, 0)
    5  = @(1, 0.4741, 1, 1, 0, 0, 0, 0, 0, 0)
    6  = @(0, 0, 0, 0, 1, 0.0025, 0.0291, 0.01, 0.0372, 0)
    7  = @(0, 0, 0, 0, 0.0025, 1, 0, 0, 0, 0)
    8  = @(0, 0, 0, 0, 0.0291, 0, 1, 1, 1, 0)
    9  = @(0, 0, 0, 0, 0.01, 0, 1, 1, 1, 0)
    10 = @(0, 0, 0, 0, 0.0372, 0, 1, 1, 1, 0)
    11 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 1)
}

# ---------------------------------------------------------------------------
# 4) Insert "sharpe_period" right before "VaR", and re-create "VaR" after it
#    with a fresh matrix of values.
#
#    The old "VaR" sheet's numbers are actually the basis for the new
#    "sharpe_period" sheet, so rather than inventing a new sheet from
#    scratch (and having to hand-roll its header row / label column / cell
#    styles to match the rest of the workbook), we rename the existing
#    "VaR" sheet in place to "sharpe_period" and refresh its values. A new,
#    blank sheet is then added right after it, renamed to "VaR", and given
#    the same header/label layout (cloned from "sharpe_period") before it
#    gets its own fresh values.
# ---------------------------------------------------------------------------
$wsVaR = $wb.Worksheets.Item("VaR")
$wsVaR.Name = "sharpe_period"
$wsSharpePeriod = $wsVaR

Set-MatrixValues $wsSharpePeriod @{
    2  = @(1, 0.8463000000000001, 1, 1, 0, 0, 0, 0, 0, 0)
    3  = @(0.8463000000000001, 1, 1, 1, 0, 0, 0, 0, 0, 0)
    4  = @(1, 1, 1, 1, 0, 0, 0, 0, 0, 0)
    5  = @(1, 1, 1, 1, 0, 0, 0, 0, 0, 0)
    6  = @(0, 0, 0, 0, 1, 0, 1, 1, 1, 0)
    7  = @(0, 0, 0, 0, 0, 1, 0, 0, 0.0001, 0)
    8  = @(0, 0, 0, 0, 1, 0, 1, 1, 1, 0)
    9  = @(0, 0, 0, 0, 1, 0, 1, 1, 1, 0)
    10 = @(0, 0, 0, 0, 1, 0.0001, 1, 1, 1, 0)
    11 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 1)
}

# New "VaR" sheet, inserted right after "sharpe_period".
$wsNewVaR = $wb.Worksheets.Add($null, $wsSharpePeriod)
$wsNewVaR.Name = "VaR"

# Clone the header row / label column (and their bold+bordered style) from
# "sharpe_period" so the new sheet matches the rest of the workbook's
# template, then clear the stray top-left corner cell picked up by the copy.
$wsSharpePeriod.Range("A1:K11").Copy($wsNewVaR.Range("A1:K11"))
$wsNewVaR.Range("A1").Clear()

Set-MatrixValues $wsNewVaR @{
    2  = @(1, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    3  = @(0, 1, 0, 0, 1, 0, 0, 0, 0, 0)
    4  = @(0, 0, 1, 1, 0.0011, 0, 0, 0, 0, 0)
    5  = @(0, 0, 1, 1, 0.001, 0, 0, 0, 0, 0)
    6  = @(0, 1, 0.0011, 0.001, 1, 0, 0, 0, 0, 0)
    7  = @(0, 0, 0, 0, 0, 1, 0, 0, 0, 0)
    8  = @(0, 0, 0, 0, 0, 0, 1, 1, 1, 0.0031)
    9  = @(0, 0, 0, 0, 0, 0, 1, 1, 1, 0.0005)
    10 = @(0, 0, 0, 0, 0, 0, 1, 1, 1, 0.0026)
    11 = @(0, 0, 0, 0, 0, 0, 0.0031, 0.0005, 0.0026, 1)
}

# Restore the originally-active tab (the sheet inserts above left the new
# "VaR" sheet selected).
$wsAnn.Select()
